# jxls sxssf_template.xlsx: add a "Total" column (E) with an =C*(1+D) formula,
# driven by a new jx:updateCell area, and retitle the sheet.
# See commit message: "SXSSF formulas support demo (via UpdateCellCommand)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E, mirroring the formatting of column D ------------------

# E1: blank banner cell -> copy D1's fill/style (style index 2) onto E1.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteAll

# E3: new "Total" header -> copy D3's header style (index 1), then set text.
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "Total"

# Retitle the banner (A1). Done *after* writing E3 so the shared-string slot
# vacated by the old title text is reused by "Total" and the new title text
# lands in a freshly appended slot (matches the workbook's string table).
$ws.Range("A1").Value = "Simple SXSSF Output"

# E4: the actual total formula, formatted with a plain 2-decimal number format.
$ws.Range("E4").Formula = "=C4*(1+D4)"
$ws.Range("E4").NumberFormat = "0.00"

# Column width / row height touch-ups for the new column/row content.
$ws.Columns("E").ColumnWidth = 13
$ws.Rows(4).RowHeight = 14.45

# --- jxls directive comments ------------------------------------------------

# Existing area/each comments now cover the extra column: lastCell D4 -> E4.
$cmtArea = $ws.Range("A1").Comment
$cmtArea.Text("Автор:
jx:area(lastCell=`"E4`")") | Out-Null

$cmtEach = $ws.Range("A4").Comment
$cmtEach.Text("Автор:
jx:each(items=`"employees`" var=`"employee`" lastCell=`"E4`")") | Out-Null

# New jx:updateCell comment driving the total column via UpdateCellCommand.
$ws.Range("E4").AddComment("Автор:
jx:updateCell(lastCell=`"E4`"  updater=`"totalCellUpdater`")") | Out-Null

# Leave selection on A1 (matches the saved workbook no longer pinning D4).
$ws.Range("A1").Select() | Out-Null
